$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns F (TVN) and G (CTC), copying header style from E1
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "TVN"
$ws.Cells.Item(1, 7).Value = "CTC"

# Rebuild data rows 2-51 with the new top-50 (time-filtered) snippet set
$data = @(
    @(2, "5:26 AM", 44.352, 74.352, 3.348, 0, 0),
    @(79, "6:05 AM", 2369.39, 2399.39, 4.37, 0, 0),
    @(103, "6:16 AM", 3079.71, 3109.71, 14.55, 0, 1),
    @(111, "6:20 AM", 3312.795, 3342.795, 0.808333, 0, 0),
    @(121, "6:25 AM", 3614.45, 3644.45, 3.595, 0, 0),
    @(132, "6:31 AM", 3944.838333, 3974.838333, 1.058333, 0, 0.166667),
    @(146, "6:38 AM", 4379.04, 4409.04, 3.49, 0, 0),
    @(169, "6:49 AM", 5046.976667, 5076.976667, 4.816667, 0, 0.333333),
    @(186, "6:58 AM", 5567.224, 5597.224, 4.226, 0, 0.2),
    @(222, "7:16 AM", 6648.33, 6678.33, 2.5925, 0, 0),
    @(234, "7:22 AM", 7006.473333, 7036.473333, 3.553333, 0, 0.333333),
    @(243, "7:26 AM", 7276.68, 7306.68, 9.06, 0, 0),
    @(251, "7:30 AM", 7516.66, 7546.66, 27.63, 0, 1),
    @(284, "7:47 AM", 8510.64, 8540.64, 0.98, 0, 0),
    @(332, "8:11 AM", 9944.236667, 9974.236667, 1.893333, 0, 0),
    @(348, "8:19 AM", 10436.51, 10466.51, 3.915, 0, 0),
    @(365, "8:27 AM", 10930.83, 10960.83, 15.57, 0, 0),
    @(378, "8:34 AM", 11323.225, 11353.225, 8.6, 0, 0),
    @(407, "8:48 AM", 12193.74, 12223.74, 1.355, 58.128333, 0),
    @(527, "9:48 AM", 15787.9, 15817.9, 2.94, 0, 0),
    @(567, "10:08 A", 16994.61, 17024.61, 7.28, 0, 0),
    @(980, "1:35 PM", 29392.66, 29422.66, 6.8, 0, 0),
    @(991, "1:40 PM", 29712.875, 29742.875, 10.58, 0, 0.5),
    @(1003, "1:46 PM", 30076.196667, 30106.196667, 3.146667, 0, 0),
    @(1030, "2:00 PM", 30884.16, 30914.16, 9.41, 0, 0.5),
    @(1049, "2:09 PM", 31454.845, 31484.845, 4.085, 0, 0),
    @(1090, "2:30 PM", 32688.42, 32718.42, 19.165, 0, 0),
    @(1104, "2:37 PM", 33100.574, 33130.574, 2.998, 0, 0),
    @(1158, "3:04 PM", 34727.39, 34757.39, 1.79, 0, 0),
    @(1176, "3:13 PM", 35251, 35281, 2.875, 0, 0),
    @(1193, "3:21 PM", 35773.465, 35803.465, 9.34, 0, 0),
    @(1202, "3:26 PM", 36041.01, 36071.01, 24.126667, 0, 0.333333),
    @(1214, "3:32 PM", 36391.61, 36421.61, 96.52, 0, 7),
    @(1230, "3:40 PM", 36884.13, 36914.13, 3.52, 0, 0.666667),
    @(1243, "3:46 PM", 37267.225, 37297.225, 9.2375, 0, 0),
    @(1276, "4:03 PM", 38263.303333, 38293.303333, 3.406667, 0, 0),
    @(1301, "4:15 PM", 39012.535, 39042.535, 7.56, 0.5, 0),
    @(1311, "4:20 PM", 39315.5, 39345.5, 76.02, 0, 1.5),
    @(1337, "4:33 PM", 40096.696667, 40126.696667, 2.046667, 0, 0),
    @(1353, "4:42 PM", 40589.99, 40619.99, 4.08, 0, 0),
    @(1414, "5:12 PM", 42400.085, 42430.085, 51.73, 0, 1),
    @(1425, "5:17 PM", 42732.5175, 42762.5175, 0.9325, 0, 0.25),
    @(1444, "5:27 PM", 43298.77, 43328.77, 16.795, 0, 0),
    @(1479, "5:44 PM", 44340.68, 44370.68, 49.65, 0, 3),
    @(1492, "5:51 PM", 44738.87, 44768.87, 22.62, 0, 0),
    @(1509, "5:59 PM", 45253.9825, 45283.9825, 3.4775, 0, 0.25),
    @(1518, "6:04 PM", 45519.58, 45549.58, 5.9, 0, 0.333333),
    @(1534, "6:12 PM", 46012.16, 46042.16, 2.17, 0, 0),
    @(1548, "6:19 PM", 46421.78, 46451.78, 0.6875, 0, 0),
    @(1569, "6:30 PM", 47068.13, 47098.13, 2.06, 0, 0),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
